# "Implement case properties (Excel, code and output)"
#
# Renames the four worksheets and retitles the generic "property" header
# column (on the case_property and step_property sheets) to "value", and
# retitles the "outliers" flag column on the step sheet to "has_outliers".

$wb = $excel.ActiveWorkbook

# --- case_property (todo) -> case_properties -------------------------------
$wsCaseProp = $wb.Worksheets.Item("case_property (todo)")
$wsCaseProp.Name = "case_properties"
$wsCaseProp.Range("B1").Value = "value"

# --- step_property (todo) -> step_properties (todo) -------------------------
$wsStepProp = $wb.Worksheets.Item("step_property (todo)")
$wsStepProp.Name = "step_properties (todo)"
$wsStepProp.Range("C1").Value = "value"

# --- step -> steps ------------------------------------------------------
$wsStep = $wb.Worksheets.Item("step")
$wsStep.Name = "steps"
$wsStep.Range("D1").Value = "has_outliers"

# --- process_flow: update the selected range left over from editing --------
$wsFlow = $wb.Worksheets.Item("process_flow")
$wsFlow.Activate()
$wsFlow.Range("A8:A9").Select()

# Return focus to the primary sheet, matching the saved workbook state.
$wsStep.Activate()

Write-Host "Renamed sheets and updated property/outlier headers."
